$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (BallBounced / Brick hit / Brick destroyed) ---
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "BallBounced"
$ws.Range("D4").Value = "Server"

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "Brick hit"
$ws.Range("D5").Value = "Server"

$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Brick destroyed"
$ws.Range("D6").Value = "Server"

# Merge SuperOp=2 rows (Brick hit / Brick destroyed), same as existing A2:A3 merge
$ws.Range("A5:A6").Merge()
$ws.Range("A5:A6").VerticalAlignment = -4108

# --- Trailing blank separator row ---
$ws.Range("A7:D7").Value = 0
$ws.Range("A7:D7").ClearContents()

# --- Fill colours: 20% - Accent6 (green, theme 9 tint 0.8 = E2EFDA) for the BallBounced row ---
$ws.Range("A4:D4").Interior.Color = 14348258

# --- Fill colours: 20% - Accent2 (orange, theme 5 tint 0.8 = FCE4D6) for the brick rows ---
$ws.Range("A2:D3").Interior.Color = 14083324
$ws.Range("A5:D6").Interior.Color = 14083324

# --- Borders: horizontal banding between blocks ---
$ws.Range("A2:D2").Borders.Item(8).LineStyle = 1
$ws.Range("A3:D3").Borders.Item(9).LineStyle = 1
$ws.Range("A4:D4").Borders.Item(8).LineStyle = 1
$ws.Range("A4:D4").Borders.Item(9).LineStyle = 1
$ws.Range("A5:D5").Borders.Item(8).LineStyle = 1
$ws.Range("D6").Borders.Item(9).LineStyle = 1
$ws.Range("A7:D7").Borders.Item(8).LineStyle = 1

# --- Borders: right edge of the Sender column (D1:D6) outlining the coloured block ---
$ws.Range("D1:D6").Borders.Item(10).LineStyle = 1

Write-Output "done"
